$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting existing rows 11:129 down to 12:130.
$ws.Rows("11:11").Insert()

# Copy the date-column number format (style) from the row above into the new row's D cell.
$ws.Range("D10").Copy()
$ws.Range("D11").PasteSpecial(-4122) | Out-Null

# Populate the new row 11 with the new record's data.
$ws.Range("A11").Value = 5
$ws.Range("B11").Value = "Macroferia Regional de Talca"
$ws.Range("C11").Value = "Maule"
$ws.Range("D11").Value = [DateTime]"2021-10-21"
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 100112017
$ws.Range("G11").Value = "Apio"
$ws.Range("H11").Value = "Americana (o)"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 500
$ws.Range("K11").Value = 7000
$ws.Range("L11").Value = 7000
$ws.Range("M11").Value = 7000
$ws.Range("N11").Value = "`$/docena de matas"
$ws.Range("O11").Value = "Provincia del Elquí"
$ws.Range("P11").Value = 1167
$ws.Range("Q11").Value = 6
$ws.Range("R11").Value = "Hortaliza"
